$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "author"
$ws.Range("I1").Value = "editor"
$ws.Range("J1").Value = "not found"
$ws.Range("J3").Value = "yes"
$ws.Range("H6").Value = "yes"
$ws.Range("H7").Value = "yes"
$ws.Range("H8").Value = "yes"
$ws.Range("H9").Value = "yes"
$ws.Range("H10").Value = "yes"
$ws.Range("H11").Value = "yes"
$ws.Range("H12").Value = "yes"
$ws.Range("H23").Value = "yes"
$ws.Range("H24").Value = "yes"
$ws.Range("H25").Value = "yes"
$ws.Range("H26").Value = "yes"
$ws.Range("H27").Value = "yes"
$ws.Range("H28").Value = "yes"
$ws.Range("H29").Value = "yes"
$ws.Range("H30").Value = "yes"
$ws.Range("H33").Value = "yes"
$ws.Range("H34").Value = "yes"
$ws.Range("H42").Value = "yes"
$ws.Range("H43").Value = "yes"
$ws.Range("H44").Value = "yes"
$ws.Range("H45").Value = "yes"
$ws.Range("H46").Value = "yes"
$ws.Range("H47").Value = "yes"
$ws.Range("H48").Value = "yes"
$ws.Range("H49").Value = "yes"
$ws.Range("H50").Value = "yes"
$ws.Range("H53").Value = "yes"
$ws.Range("H54").Value = "yes"
$ws.Range("H86").Value = "yes"
$ws.Range("H89").Value = "yes"
$ws.Range("H90").Value = "yes"
$ws.Range("H91").Value = "yes"
$ws.Range("H92").Value = "yes"
$ws.Range("H93").Value = "yes"
$ws.Range("H94").Value = "yes"
$ws.Range("H95").Value = "yes"
$ws.Range("H96").Value = "yes"
$ws.Range("H97").Value = "yes"
$ws.Range("H98").Value = "yes"
$ws.Range("H99").Value = "yes"
$ws.Range("H116").Value = "yes"
$ws.Range("H118").Value = "yes"
$ws.Range("H119").Value = "yes"
$ws.Range("H120").Value = "yes"
$ws.Range("H121").Value = "yes"
$ws.Range("H122").Value = "yes"
$ws.Range("H123").Value = "yes"
$ws.Range("H124").Value = "yes"
$ws.Range("H157").Value = "yes"
$ws.Range("H158").Value = "yes"
$ws.Range("H159").Value = "yes"
$ws.Range("H160").Value = "yes"
$ws.Range("H161").Value = "yes"
$ws.Range("H162").Value = "yes"
$ws.Range("H163").Value = "yes"
$ws.Range("H164").Value = "yes"
$ws.Range("H188").Value = "yes"
$ws.Range("I198").Value = "yes"
$ws.Range("H210").Value = "yes"
$ws.Range("H232").Value = "yes"
$ws.Range("H233").Value = "yes"
$ws.Range("H236").Value = "yes"
$ws.Range("H237").Value = "yes"
$ws.Range("H238").Value = "yes"
$ws.Range("H239").Value = "yes"
$ws.Range("H240").Value = "yes"
$ws.Range("H241").Value = "yes"
$ws.Range("H242").Value = "yes"
$ws.Range("H243").Value = "yes"
$ws.Range("H256").Value = "yes"
$ws.Range("H259").Value = "yes"
$ws.Range("H260").Value = "yes"
$ws.Range("H261").Value = "yes"
$ws.Range("H262").Value = "yes"
$ws.Range("H263").Value = "yes"
$ws.Range("H264").Value = "yes"
$ws.Range("H265").Value = "yes"
$ws.Range("H279").Value = "yes"
$ws.Range("H295").Value = "yes"

# Apply green (theme accent6) font color to flagged D cells
$ws.Range("D53").Font.ThemeColor = 10
$ws.Range("D54").Font.ThemeColor = 10
$ws.Range("D116").Font.ThemeColor = 10
$ws.Range("D188").Font.ThemeColor = 10
$ws.Range("D198").Font.ThemeColor = 10
$ws.Range("D233").Font.ThemeColor = 10
$ws.Range("D258").Font.ThemeColor = 10
